$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column L: header + a date value copied (format-wise) from D2 (DATE_START_OF_COVER)
$ws.Range("L1").Value = "DATE_OF_DISABLEMENT"
$ws.Range("L2").Value = $ws.Range("D2").Value2

$ws.Range("D2").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("L2").Select()
